# fix: fixed formatting when scrapping floating point numbers
#
# The "Importe" column (H) holds amounts that were scraped with
# Spanish/Argentine thousands ('.') and decimal (',') separators, e.g.
# "28.248,75". This normalizes them to plain decimal-point numeric text,
# e.g. "28248.75" - still stored as TEXT (not converted to a real number),
# matching the original authoring.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> corrected text value (old Spanish-formatted -> new).
$fixes = [ordered]@{
    "H2"  = "28248.75"
    "H3"  = "216000.00"
    "H4"  = "130380.00"
    "H5"  = "450.00"
    "H6"  = "9952.00"
    "H7"  = "75150.00"
    "H8"  = "13800.00"
    "H9"  = "31974.00"
    "H10" = "13800.00"
    "H11" = "50000.00"
    "H12" = "21000.00"
    "H13" = "61320.00"
    "H14" = "12094.06"
    "H15" = "213000.00"
    "H16" = "1842.68"
    "H17" = "6530.00"
}

foreach ($addr in $fixes.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    # Leading apostrophe forces Excel to keep the numeric-looking text as a
    # literal string (t="s" / shared string) instead of coercing it to a
    # real number.
    $cell.Value = "'" + $fixes[$addr]
    # Restore the original (default) style so we don't leave a stray
    # "quote prefix" number format applied to the cell.
    $cell.Style = $origStyle
}
